$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.1425304360311941
$ws.Range("D2").Value = 0.00136281687266937
$ws.Range("E2").Value = 0.4336984297763706
$ws.Range("F2").Value = 0.7160575709075232
$ws.Range("G2").Value = 0.002324725457951165
$ws.Range("N2").Value = 1.995241137946778
$ws.Range("O2").Value = 2.355210718842727

$ws.Range("B3").Value = 0.1329091933101978
$ws.Range("D3").Value = 0.001187382742090293
$ws.Range("E3").Value = 0.3778594790115903
$ws.Range("F3").Value = 0.6567163053088905
$ws.Range("G3").Value = 0.002329461130050568
$ws.Range("N3").Value = 1.846526605564804
$ws.Range("O3").Value = 2.156718067768907

$ws.Range("B4").Value = 0.1270661492628875
$ws.Range("D4").Value = 0.001080080274057593
$ws.Range("E4").Value = 0.3437100912600357
$ws.Range("F4").Value = 0.6207120828670867
$ws.Range("G4").Value = 0.002332519531335734
$ws.Range("N4").Value = 1.755520437593759
$ws.Range("O4").Value = 2.036265997455075

$ws.Range("B5").Value = 0.1247015078381395
$ws.Range("D5").Value = 0.001036432189167158
$ws.Range("E5").Value = 0.329825076174572
$ws.Range("F5").Value = 0.6061468742205847
$ws.Range("G5").Value = 0.002333803886885216
$ws.Range("N5").Value = 1.718516884911963
$ws.Range("O5").Value = 1.987532592803802

$ws.Range("B6").Value = 0.1243098614095572
$ws.Range("D6").Value = 0.001029188458627317
$ws.Range("E6").Value = 0.3275212721868144
$ws.Range("F6").Value = 0.6037347334826251
$ws.Range("G6").Value = 0.002334019454337199
$ws.Range("N6").Value = 1.712377599108464
$ws.Range("O6").Value = 1.979461519004303

$ws.Range("B7").Value = 0.1270341919673399
$ws.Range("D7").Value = 0.001079491335964633
$ws.Range("E7").Value = 0.3435227109891912
$ws.Range("F7").Value = 0.6205152211821883
$ws.Range("G7").Value = 0.002332536698358234
$ws.Range("N7").Value = 1.755021054648495
$ws.Range("O7").Value = 2.035607345351195

$ws.Range("B8").Value = 0.1391997998654944
$ws.Range("D8").Value = 0.001302225393729373
$ws.Range("E8").Value = 0.4144150822294819
$ws.Range("F8").Value = 0.6955060775010651
$ws.Range("G8").Value = 0.002326327130626193
$ws.Range("N8").Value = 1.943904251618477
$ws.Range("O8").Value = 2.286471316272696

$ws.Range("B9").Value = 0.1635584291339143
$ws.Range("D9").Value = 0.001743495443395915
$ws.Range("E9").Value = 0.5546583365363915
$ws.Range("F9").Value = 0.8460728912957052
$ws.Range("G9").Value = 0.002315339238254839
$ws.Range("N9").Value = 2.316504389646809
$ws.Range("O9").Value = 2.79001036887513

$ws.Range("B10").Value = 0.1817500012197257
$ws.Range("D10").Value = 0.002072145766488731
$ws.Range("E10").Value = 0.658652049917535
$ws.Range("F10").Value = 0.9589651656520459
$ws.Range("G10").Value = 0.002307982165642218
$ws.Range("N10").Value = 2.591327722520418
$ws.Range("O10").Value = 3.167490004210322

$ws.Range("B11").Value = 0.1900877858899008
$ws.Range("D11").Value = 0.002223006803760086
$ws.Range("E11").Value = 0.7062161564137739
$ws.Range("F11").Value = 1.010845266279489
$ws.Range("G11").Value = 0.002304788697644411
$ws.Range("N11").Value = 2.716530669508643
$ws.Range("O11").Value = 3.340953389518518

$ws.Range("B12").Value = 0.1932538308938092
$ws.Range("D12").Value = 0.002280360919026947
$ws.Range("E12").Value = 0.7242680764979781
$ws.Range("F12").Value = 1.03056856590527
$ws.Range("G12").Value = 0.00230360130780261
$ws.Range("N12").Value = 2.763963193734583
$ws.Range("O12").Value = 3.406898250946483

$ws.Range("B13").Value = 0.1925715833279469
$ws.Range("D13").Value = 0.002267998111934233
$ws.Range("E13").Value = 0.720378418342051
$ws.Range("F13").Value = 1.026317325067765
$ws.Range("G13").Value = 0.00230385606152078
$ws.Range("N13").Value = 2.753746911385861
$ws.Range("O13").Value = 3.392684256368284

$ws.Range("B14").Value = 0.190348085298524
$ws.Range("D14").Value = 0.00222772064371668
$ws.Range("E14").Value = 0.7077004683360144
$ws.Range("F14").Value = 1.012466350817789
$ws.Range("G14").Value = 0.002304690571907265
$ws.Range("N14").Value = 2.720432585385993
$ws.Range("O14").Value = 3.346373501059304

$ws.Range("B15").Value = 0.1889872547792919
$ws.Range("D15").Value = 0.002203079946879072
$ws.Range("E15").Value = 0.6999402271471666
$ws.Range("F15").Value = 1.003992371776604
$ws.Range("G15").Value = 0.002305204582768017
$ws.Range("N15").Value = 2.700029152110062
$ws.Range("O15").Value = 3.318040639829917

$ws.Range("B16").Value = 0.1812063398625838
$ws.Range("D16").Value = 0.002062316590961899
$ws.Range("E16").Value = 0.655549089093725
$ws.Range("F16").Value = 0.9555854069586189
$ws.Range("G16").Value = 0.002308193939223566
$ws.Range("N16").Value = 2.583148644061055
$ws.Range("O16").Value = 3.15618949106215

$ws.Range("B17").Value = 0.1764487851811651
$ws.Range("D17").Value = 0.001976332213386556
$ws.Range("E17").Value = 0.6283848123191405
$ws.Range("F17").Value = 0.9260251098173171
$ws.Range("G17").Value = 0.002310066979291825
$ws.Range("N17").Value = 2.511489453097226
$ws.Range("O17").Value = 3.057351167739682

$ws.Range("B18").Value = 0.1737182494248515
$ws.Range("D18").Value = 0.001927002019733948
$ws.Range("E18").Value = 0.6127846016899383
$ws.Range("F18").Value = 0.9090720897889923
$ws.Range("G18").Value = 0.00231115873986465
$ws.Range("N18").Value = 2.470290675909951
$ws.Range("O18").Value = 3.000665896692624

$ws.Range("B19").Value = 0.1727947553757332
$ws.Range("D19").Value = 0.001910320321501402
$ws.Range("E19").Value = 0.6075066487272522
$ws.Range("F19").Value = 0.9033404933890949
$ws.Range("G19").Value = 0.002311530874994851
$ws.Range("N19").Value = 2.456344689235948
$ws.Range("O19").Value = 2.981501180670932

$ws.Range("B20").Value = 0.1769546279144265
$ws.Range("D20").Value = 0.001985472159233836
$ws.Range("E20").Value = 0.6312739907922804
$ws.Range("F20").Value = 0.9291667394193155
$ws.Range("G20").Value = 0.002309866097720437
$ws.Range("N20").Value = 2.519115905984052
$ws.Range("O20").Value = 3.067855660969485

$ws.Range("B21").Value = 0.1910009463252038
$ws.Range("D21").Value = 0.002239544716820419
$ws.Range("E21").Value = 0.7114231643078455
$ws.Range("F21").Value = 1.016532600446681
$ws.Range("G21").Value = 0.002304444862530205
$ws.Range("N21").Value = 2.730217288371648
$ws.Range("O21").Value = 3.359969033151685

$ws.Range("B22").Value = 0.2002316580432932
$ws.Range("D22").Value = 0.002406930069309254
$ws.Range("E22").Value = 0.7640425772889614
$ws.Range("F22").Value = 1.074083505283141
$ws.Range("G22").Value = 0.002301029401395692
$ws.Range("N22").Value = 2.868303394937072
$ws.Range("O22").Value = 3.552389396738079

$ws.Range("B23").Value = 0.195300506148314
$ws.Range("D23").Value = 0.00231746072000405
$ws.Range("E23").Value = 0.7359357237643849
$ws.Range("F23").Value = 1.043325465387227
$ws.Range("G23").Value = 0.002302840664474282
$ws.Range("N23").Value = 2.794595190493112
$ws.Range("O23").Value = 3.449550766766095

$ws.Range("B24").Value = 0.1767259219242163
$ws.Range("D24").Value = 0.001981339670157567
$ws.Range("E24").Value = 0.6299677408469222
$ws.Range("F24").Value = 0.927746279284861
$ws.Range("G24").Value = 0.002309956869657329
$ws.Range("N24").Value = 2.51566798917645
$ws.Range("O24").Value = 3.06310614965065

$ws.Range("B25").Value = 0.1569161479393699
$ws.Range("D25").Value = 0.00162346534227531
$ws.Range("E25").Value = 0.5165662124021679
$ws.Range("F25").Value = 0.804951376280016
$ws.Range("G25").Value = 0.002318185404700003
$ws.Range("N25").Value = 1.995241137946778
$ws.Range("O25").Value = 2.652501122830984
